$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(2833958, 62.9, 1783065, 31.6, 896590, 4.6, 129949, 0.9, 24355, 37.1, 1050893)
    3 = @(2471589, 66.5, 1644768, 28.5, 703296, 4.3, 107447, 0.7, 16079, 33.5, 826821)
    4 = @(226810, 36.8, 83383, 54.8, 124350, 5.4, 12235, 3, 6842, 63.2, 143427)
    5 = @(113405, 42.2, 47820, 48.9, 55452, 7.7, 8760, 1.2, 1373, 57.8, 65585)
    6 = @(22154, 32, 7094, 60.9, 13492, 6.8, 1507, 0.3, 62, 68, 15060)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 3]
    }
}

$wb.Save()
